$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed values (per row, columns D:AJ of the IFRS consolidated financial table)

# Row 2
$ws.Range("D2").Value = 148454
$ws.Range("E2").Value = 4844
$ws.Range("F2").Value = 4844
$ws.Range("G2").Value = 5244
$ws.Range("H2").Value = 4051
$ws.Range("I2").Value = 4054
$ws.Range("J2").Value = -3
$ws.Range("K2").Value = 921262
$ws.Range("L2").Value = 834741
$ws.Range("M2").Value = 86522
$ws.Range("N2").Value = 86516
$ws.Range("O2").Value = 6
$ws.Range("P2").Value = 43426
$ws.Range("Q2").Value = 29777
$ws.Range("R2").Value = -27402
$ws.Range("S2").Value = -1541
$ws.Range("T2").Value = 2652
$ws.Range("V2").Value = 1737
$ws.Range("W2").Value = 3.26
$ws.Range("X2").Value = 2.73
$ws.Range("Y2").Value = 5.18
$ws.Range("Z2").Value = 0.46
$ws.Range("AA2").Value = 964.78
$ws.Range("AB2").Value = 106.98
$ws.Range("AC2").Value = 467
$ws.Range("AD2").Value = 17.76
$ws.Range("AE2").Value = 10487
$ws.Range("AF2").Value = 0.79
$ws.Range("AG2").Value = 180
$ws.Range("AH2").Value = 2.17
$ws.Range("AI2").Value = 36.71
$ws.Range("AJ2").Value = 868530000

# Row 3
$ws.Range("D3").Value = 169409
$ws.Range("E3").Value = 5866
$ws.Range("F3").Value = 5866
$ws.Range("G3").Value = 6070
$ws.Range("H3").Value = 5300
$ws.Range("I3").Value = 5301
$ws.Range("J3").Value = -1
$ws.Range("K3").Value = 995302
$ws.Range("L3").Value = 908874
$ws.Range("M3").Value = 86428
$ws.Range("N3").Value = 86423
$ws.Range("O3").Value = 5
$ws.Range("P3").Value = 43426
$ws.Range("Q3").Value = 82070
$ws.Range("R3").Value = -71297
$ws.Range("S3").Value = -7156
$ws.Range("T3").Value = 1701
$ws.Range("V3").Value = 1983
$ws.Range("W3").Value = 3.46
$ws.Range("X3").Value = 3.13
$ws.Range("Y3").Value = 6.13
$ws.Range("Z3").Value = 0.55
$ws.Range("AA3").Value = 1051.6
$ws.Range("AB3").Value = 120.3
$ws.Range("AC3").Value = 610
$ws.Range("AD3").Value = 12.11
$ws.Range("AE3").Value = 11502
$ws.Range("AF3").Value = 0.64
$ws.Range("AG3").Value = 180
$ws.Range("AH3").Value = 2.44
$ws.Range("AI3").Value = 25.51
$ws.Range("AJ3").Value = 868530000

# Row 4
$ws.Range("D4").Value = 206661
$ws.Range("E4").Value = 5210
$ws.Range("F4").Value = 5210
$ws.Range("G4").Value = 10093
$ws.Range("H4").Value = 8451
$ws.Range("I4").Value = 7953
$ws.Range("J4").Value = 498
$ws.Range("K4").Value = 1193811
$ws.Range("L4").Value = 1102220
$ws.Range("M4").Value = 91592
$ws.Range("N4").Value = 85291
$ws.Range("O4").Value = 6301
$ws.Range("P4").Value = 43426
$ws.Range("Q4").Value = 40878
$ws.Range("R4").Value = -38236
$ws.Range("S4").Value = -1738
$ws.Range("T4").Value = 690
$ws.Range("V4").Value = 4350
$ws.Range("W4").Value = 2.52
$ws.Range("X4").Value = 4.09
$ws.Range("Y4").Value = 9.84
$ws.Range("Z4").Value = 0.77
$ws.Range("AA4").Value = 1203.41
$ws.Range("AB4").Value = 132.19
$ws.Range("AC4").Value = 916
$ws.Range("AD4").Value = 7.13
$ws.Range("AE4").Value = 11351
$ws.Range("AF4").Value = 0.58
$ws.Range("AG4").Value = 80
$ws.Range("AH4").Value = 1.23
$ws.Range("AI4").Value = 7.56
$ws.Range("AJ4").Value = 868530000

# Row 5
$ws.Range("D5").Value = 260871
$ws.Range("E5").Value = 9534
$ws.Range("F5").Value = 9534
$ws.Range("G5").Value = 9645
$ws.Range("H5").Value = 6887
$ws.Range("I5").Value = 5849
$ws.Range("J5").Value = 1038
$ws.Range("K5").Value = 1259945
$ws.Range("L5").Value = 1157564
$ws.Range("M5").Value = 102381
$ws.Range("N5").Value = 93516
$ws.Range("O5").Value = 8865
$ws.Range("P5").Value = 43426
$ws.Range("Q5").Value = 33187
$ws.Range("R5").Value = -40662
$ws.Range("S5").Value = 4806
$ws.Range("T5").Value = 1137
$ws.Range("V5").Value = 3419
$ws.Range("W5").Value = 3.65
$ws.Range("X5").Value = 2.64
$ws.Range("Y5").Value = 7.7
$ws.Range("Z5").Value = 0.5600000000000001
$ws.Range("AA5").Value = 1130.64
$ws.Range("AB5").Value = 157.04
$ws.Range("AC5").Value = 673
$ws.Range("AD5").Value = 10.26
$ws.Range("AE5").Value = 12446
$ws.Range("AF5").Value = 0.5600000000000001
$ws.Range("AG5").Value = 140
$ws.Range("AH5").Value = 2.03
$ws.Range("AI5").Value = 17.98
$ws.Range("AJ5").Value = 868530000

# Row 6
$ws.Range("D6").Value = 234305
$ws.Range("E6").Value = 6502
$ws.Range("F6").Value = 6502
$ws.Range("G6").Value = 6364
$ws.Range("H6").Value = 4465
$ws.Range("I6").Value = 4153
$ws.Range("K6").Value = 1320845
$ws.Range("L6").Value = 1206548
$ws.Range("M6").Value = 114297
$ws.Range("N6").Value = 103069
$ws.Range("P6").Value = 43426
$ws.Range("Q6").Value = 16429
$ws.Range("R6").Value = -18322
$ws.Range("S6").Value = 14334
$ws.Range("T6").Value = 1673
$ws.Range("V6").Value = 6903
$ws.Range("W6").Value = 2.77
$ws.Range("X6").Value = 1.91
$ws.Range("Y6").Value = 4.54
$ws.Range("Z6").Value = 0.35
$ws.Range("AA6").Value = 1055.63
$ws.Range("AB6").Value = 184.47
$ws.Range("AC6").Value = 478
$ws.Range("AD6").Value = 8.83
$ws.Range("AE6").Value = 13717
$ws.Range("AF6").Value = 0.31
$ws.Range("AG6").Value = 100
$ws.Range("AH6").Value = 2.37
$ws.Range("AI6").Value = 18.09
$ws.Range("AJ6").Value = 868530000

# Row 7
$ws.Range("D7").Value = 95530
$ws.Range("E7").Value = 2330
$ws.Range("G7").Value = 2090
$ws.Range("H7").Value = 1870
$ws.Range("I7").Value = 1470
$ws.Range("K7").Value = 1390155
$ws.Range("L7").Value = 1267210
$ws.Range("M7").Value = 122945
$ws.Range("N7").Value = 110825
$ws.Range("P7").Value = 43430
$ws.Range("W7").Value = 2.44
$ws.Range("X7").Value = 1.96
$ws.Range("Y7").Value = 1.38
$ws.Range("Z7").Value = 0.14
$ws.Range("AA7").Value = 1030.71
$ws.Range("AC7").Value = 169
$ws.Range("AD7").Value = 12.11
$ws.Range("AE7").Value = 14749
$ws.Range("AF7").Value = 0.14
$ws.Range("AG7").Value = 33
$ws.Range("AH7").Value = 1.63
$ws.Range("AI7").Value = 19.71

# Row 8
$ws.Range("D8").Value = 96815
$ws.Range("E8").Value = 4140
$ws.Range("G8").Value = 4895
$ws.Range("H8").Value = 3605
$ws.Range("I8").Value = 2745
$ws.Range("K8").Value = 1439845
$ws.Range("L8").Value = 1314675
$ws.Range("M8").Value = 125165
$ws.Range("N8").Value = 112905
$ws.Range("P8").Value = 43430
$ws.Range("W8").Value = 4.28
$ws.Range("X8").Value = 3.72
$ws.Range("Y8").Value = 2.45
$ws.Range("Z8").Value = 0.26
$ws.Range("AA8").Value = 1050.35
$ws.Range("AC8").Value = 316
$ws.Range("AD8").Value = 6.49
$ws.Range("AE8").Value = 15026
$ws.Range("AF8").Value = 0.14
$ws.Range("AG8").Value = 76
$ws.Range("AH8").Value = 3.73
$ws.Range("AI8").Value = 24.16

# Row 9
$ws.Range("D9").Value = 98810
$ws.Range("E9").Value = 4630
$ws.Range("G9").Value = 5485
$ws.Range("H9").Value = 4040
$ws.Range("I9").Value = 3180
$ws.Range("K9").Value = 1492030
$ws.Range("L9").Value = 1364440
$ws.Range("M9").Value = 127585
$ws.Range("N9").Value = 115180
$ws.Range("P9").Value = 43430
$ws.Range("W9").Value = 4.69
$ws.Range("X9").Value = 4.09
$ws.Range("Y9").Value = 2.79
$ws.Range("Z9").Value = 0.28
$ws.Range("AA9").Value = 1069.44
$ws.Range("AC9").Value = 366
$ws.Range("AD9").Value = 5.6
$ws.Range("AE9").Value = 15329
$ws.Range("AF9").Value = 0.13
$ws.Range("AG9").Value = 87
$ws.Range("AH9").Value = 4.22
$ws.Range("AI9").Value = 23.64

# Remove cells that no longer exist in the target layout (columns merged/dropped)
$ws.Range("U2").ClearContents()
$ws.Range("U3").ClearContents()
$ws.Range("U4").ClearContents()
$ws.Range("U5").ClearContents()
$ws.Range("U6").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
